# Add "EndUserPO" column (AC) with a value, add a CustomerPO value for row 3 (H3),
# and update the IMOrderNo value in C3, per commit "end user po column added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing IMOrderNo value on row 3.
$ws.Range("C3").Value = "37-26283-11"

# Populate the CustomerPO cell for row 3 (previously empty).
$ws.Range("H3").Value = "abc"

# Add the new "EndUserPO" column header in AC1, copying the style used by the
# other header cells (e.g. AA1) so formatting matches.
$ws.Range("AA1").Copy()
$ws.Range("AC1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("AC1").Value = "EndUserPO"

# Add the new EndUserPO value for row 3.
$ws.Range("AC3").Value = "EPO123"

# Update the sheet view so the new column is visible/selected, matching the
# saved view state from the authored workbook (scrolled so column S is the
# left-most visible column, with AC3 as the active selection).
$ws.Range("AC3").Select()
$excel.ActiveWindow.ScrollColumn = 19
